$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected values (NN correction) for rows 2-9, columns B-H (G unchanged)
$data = @{
    2 = @{ B = 1214.601440429688; C = 0.9112; D = 0.8974999785423279; E = 1.259999990463257; F = 0.8033999800682068; H = 0.6331 }
    3 = @{ B = 1208.54736328125;  C = 0.9668; D = 0.9275;              E = 1.479599952697754; F = 0.781499981880188;  H = 0.8986 }
    4 = @{ B = 817.7930297851562; C = 0.9554; D = 0.9193;              E = 1.36679995059967;  F = 0.8328999876976013; H = 0.8262 }
    5 = @{ B = 894.7968139648438; C = 0.9439; D = 0.9292;              E = 1.174299955368042; F = 0.7394000291824341; H = 0.9141 }
    6 = @{ B = 1181.625854492188; C = 0.9348; D = 0.9362;              E = 1.115599989891052; F = 0.8299999833106995; H = 0.976 }
    7 = @{ B = 906.6458740234375; C = 0.913;  D = 0.9093999862670898;  E = 1.077900052070618; F = 0.7694000005722046; H = 0.7385 }
    8 = @{ B = 993.2744140625;    C = 0.89;   D = 0.8869;              E = 1.068600058555603; F = 0.789900004863739;  H = 0.5394 }
    9 = @{ B = 7217.28515625;     C = 0.9301; D = 0.9121;              E = 1.479599952697754; F = 0.7394000291824341; H = 5.525899999999999 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
